# Applies the permutation of artfynd rows 8-24 (species records reshuffled
# across rows) plus the corresponding cell-level content updates, per the
# target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A8").Value = 111742261
$ws.Range("B8").Value = 94125
$ws.Range("E8").Value = 54
$ws.Range("F8").Value = "Skogstrappmossa"
$ws.Range("G8").Value = "Anastrophyllum michauxii"
$ws.Range("H8").Value = "(F.Weber.) H.Buch"
$ws.Range("L8").Value = ""
$ws.Range("Q8").Value = 331832.2172148526
$ws.Range("R8").Value = 6626574.972218169
$ws.Range("A9").Value = 111741735
$ws.Range("B9").Value = 73689
$ws.Range("E9").Value = 308
$ws.Range("F9").Value = "Brunpudrad nållav"
$ws.Range("G9").Value = "Chaenotheca gracillima"
$ws.Range("H9").Value = "(Vain.) Tibell"
$ws.Range("J9").Value = ""
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("Q9").Value = 331238.0752669616
$ws.Range("R9").Value = 6626585.695077832
$ws.Range("AC9").Value = "På gran"
$ws.Range("AF9").Value = ""
$ws.Range("A10").Value = 111741661
$ws.Range("B10").Value = 94134
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 53
$ws.Range("F10").Value = "Vedtrappmossa"
$ws.Range("G10").Value = "Crossocalyx hellerianus"
$ws.Range("H10").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("L10").Value = ""
$ws.Range("Q10").Value = 331248.4265637345
$ws.Range("R10").Value = 6626657.912916132
$ws.Range("AC10").Value = ""
$ws.Range("A11").Value = 111741430
$ws.Range("B11").Value = 56398
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("J11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = "färska spår"
$ws.Range("Q11").Value = 331285.2567537006
$ws.Range("R11").Value = 6626678.453820148
$ws.Range("AC11").Value = ""
$ws.Range("AF11").Value = ""
$ws.Range("A12").Value = 111741461
$ws.Range("B12").Value = 94851
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 2569
$ws.Range("F12").Value = "Stor revmossa"
$ws.Range("G12").Value = "Bazzania trilobata"
$ws.Range("H12").Value = "(L.) Gray"
$ws.Range("Q12").Value = 331282.7667812487
$ws.Range("R12").Value = 6626634.652084536
$ws.Range("A13").Value = 111741381
$ws.Range("B13").Value = 77515
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("H13").Value = "(Ach.) Ach."
$ws.Range("L13").Value = ""
$ws.Range("Q13").Value = 331269.3053187408
$ws.Range("R13").Value = 6626728.138093493
$ws.Range("A14").Value = 111741744
$ws.Range("B14").Value = 73689
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 308
$ws.Range("F14").Value = "Brunpudrad nållav"
$ws.Range("G14").Value = "Chaenotheca gracillima"
$ws.Range("H14").Value = "(Vain.) Tibell"
$ws.Range("Q14").Value = 331225.6043252015
$ws.Range("R14").Value = 6626555.973579897
$ws.Range("AC14").Value = "På björkhögstubbe"
$ws.Range("A15").Value = 111742234
$ws.Range("B15").Value = 77515
$ws.Range("E15").Value = 6425
$ws.Range("F15").Value = "Garnlav"
$ws.Range("G15").Value = "Alectoria sarmentosa"
$ws.Range("H15").Value = "(Ach.) Ach."
$ws.Range("Q15").Value = 331746.1491186697
$ws.Range("R15").Value = 6626673.250118625
$ws.Range("A16").Value = 111741405
$ws.Range("Q16").Value = 331300.7315557983
$ws.Range("R16").Value = 6626707.533622785
$ws.Range("AC16").Value = "På björkhögstubbe i fuktigt läge"
$ws.Range("A17").Value = 111741468
$ws.Range("B17").Value = 92683
$ws.Range("E17").Value = 2362
$ws.Range("F17").Value = "Blek stjärnmossa"
$ws.Range("G17").Value = "Mnium stellare"
$ws.Range("H17").Value = "Hedw."
$ws.Range("A18").Value = 111741493
$ws.Range("B18").Value = 77604
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 6450
$ws.Range("F18").Value = "Skuggblåslav"
$ws.Range("G18").Value = "Hypogymnia vittata"
$ws.Range("H18").Value = "(Ach.) Parrique"
$ws.Range("Q18").Value = 331282.7667812487
$ws.Range("R18").Value = 6626634.652084536
$ws.Range("AC18").Value = ""
$ws.Range("A19").Value = 111742250
$ws.Range("B19").Value = 77515
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6425
$ws.Range("F19").Value = "Garnlav"
$ws.Range("G19").Value = "Alectoria sarmentosa"
$ws.Range("H19").Value = "(Ach.) Ach."
$ws.Range("L19").Value = ""
$ws.Range("Q19").Value = 331780.8802231384
$ws.Range("R19").Value = 6626591.424151366
$ws.Range("A20").Value = 111741395
$ws.Range("B20").Value = 73510
$ws.Range("D20").Value = "LC"
$ws.Range("E20").Value = 6428
$ws.Range("F20").Value = "Rostfläck"
$ws.Range("G20").Value = "Arthonia vinosa"
$ws.Range("H20").Value = "Leight."
$ws.Range("Q20").Value = 331269.3053187408
$ws.Range("R20").Value = 6626728.138093493
$ws.Range("AC20").Value = "På klibbal"
$ws.Range("A21").Value = 111742256
$ws.Range("B21").Value = 79444
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 1049
$ws.Range("F21").Value = "Kortskaftad ärgspik"
$ws.Range("G21").Value = "Microcalicium ahlneri"
$ws.Range("H21").Value = "Tibell"
$ws.Range("Q21").Value = 331773.1827125447
$ws.Range("R21").Value = 6626566.53343309
$ws.Range("AC21").Value = "På barklös talltorraka"
$ws.Range("A22").Value = 111742228
$ws.Range("B22").Value = 90689
$ws.Range("E22").Value = 5966
$ws.Range("F22").Value = "Motaggsvamp"
$ws.Range("G22").Value = "Sarcodon squamosus"
$ws.Range("H22").Value = "(Schaeff.) Quél."
$ws.Range("Q22").Value = 331723.8827412428
$ws.Range("R22").Value = 6626661.637235454
$ws.Range("AC22").Value = ""
$ws.Range("A23").Value = 111741438
$ws.Range("B23").Value = 73634
$ws.Range("D23").Value = "LC"
$ws.Range("E23").Value = 6426
$ws.Range("F23").Value = "Kattfotslav"
$ws.Range("G23").Value = "Felipes leucopellaeus"
$ws.Range("H23").Value = "(Ach.) Frisch & G.Thor"
$ws.Range("Q23").Value = 331285.2567537006
$ws.Range("R23").Value = 6626678.453820148
$ws.Range("A24").Value = 111741759
$ws.Range("B24").Value = 77604
$ws.Range("D24").Value = "LC"
$ws.Range("E24").Value = 6450
$ws.Range("F24").Value = "Skuggblåslav"
$ws.Range("G24").Value = "Hypogymnia vittata"
$ws.Range("H24").Value = "(Ach.) Parrique"
$ws.Range("Q24").Value = 331225.6043252015
$ws.Range("R24").Value = 6626555.973579897
$ws.Range("AC24").Value = "På liten lodyta"
